$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet "Sheet1" -> "Sam Curran"
$ws.Name = "Sam Curran"

# 2. Insert a new column before column A for "matchNo"; existing B:L shift to C:M
$ws.Columns.Item(1).Insert()

# 3. Insert 3 new rows above the single existing data row (old row 2 becomes row 5)
$ws.Rows.Item(2).Resize(3).Insert()

# 4. Write the header row + all 4 data rows.
#    Every value is poked in as a ="literal text" formula so Excel stores it with
#    the "string" cell type (the source keeps purely-numeric fields like runs/sr as
#    text too, e.g. sr "125.00" would lose its trailing zeros as a real number),
#    then the whole block is flattened to plain literal values via a single
#    copy / paste-special pass so no formulas are left behind.
$grid = @(
    @('="matchNo"', '="teamName"', '="batterName"', '="states"', '="runs"', '="balls"', '="fours"', '="sixes"', '="sr"', '="opponentTeamName"', '="venue"', '="date"', '="result"'),
    @('="2nd"', '="Chennai Super Kings"', '="Sam Curran"', '="b Woakes"', '="34"', '="15"', '="4"', '="2"', '="226.66"', '="Delhi Capitals"', '="Wankhede"', '="April 10"', '="Capitals won by 7 wickets (with 8 balls remaining)"'),
    @('="8th"', '="Chennai Super Kings"', '="Sam Curran"', '=""', '="5"', '="4"', '="1"', '="0"', '="125.00"', '="Punjab Kings"', '="Wankhede"', '="April 16"', '="Super Kings won by 6 wickets (with 26 balls remaining)"'),
    @('="38th"', '="Chennai Super Kings"', '="Sam Curran"', '="c sub (KL Nagarkoti) b Narine"', '="4"', '="4"', '="0"', '="0"', '="100.00"', '="Kolkata Knight Riders"', '="Abu Dhabi"', '="September 26"', '="Super Kings won by 2 wickets"'),
    @('="12th"', '="Chennai Super Kings"', '="Sam Curran"', '="run out (Mustafizur Rahman/†Samson)"', '="13"', '="6"', '="0"', '="1"', '="216.66"', '="Rajasthan Royals"', '="Wankhede"', '="April 19"', '="Super Kings won by 45 runs"')
)

for ($r = 0; $r -lt $grid.Length; $r++) {
    $row = $grid[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Formula = $row[$c]
    }
}

# 5. Flatten the helper formulas down to literal values in a single bulk operation
#    (xlPasteValues = -4163).
$full = $ws.Range("A1:M5")
$full.Copy()
$full.PasteSpecial(-4163)
$excel.CutCopyMode = 0
